# Fill in the "Rui ro 3" and "Rui ro 4" risk-list paragraphs with their
# content, and strip the italic direct-formatting that the placeholder
# bullets had (the new text is regular, non-italic).
$d = $word.ActiveDocument

# --- Rui ro 3 (paragraphs 151-155) ---

$p = $d.Paragraphs.Item(151)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Tên rủi ro: không thể chạy được trên Linux")

$p = $d.Paragraphs.Item(152)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Mô tả rủi ro: không thể sử dụng chương trình trên Linux")

$p = $d.Paragraphs.Item(153)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Xác suất xảy ra: 20%")

$p = $d.Paragraphs.Item(154)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Mức độ thiệt hại: thấp")

$p = $d.Paragraphs.Item(155)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Giải pháp xử lý: Cài them thư viện để chạy chương trình")

# --- Rui ro 4 (paragraphs 157-161) ---

$p = $d.Paragraphs.Item(157)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Tên rủi ro: không thể tương tác giữa người dùng window và linux")

$p = $d.Paragraphs.Item(158)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Mô tả rủi ro: người dùng window không thể tương tác với người dùng bên linux thông qua chương trình")

$p = $d.Paragraphs.Item(159)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Xác suất xảy ra: 5%")

$p = $d.Paragraphs.Item(160)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Mức độ thiệt hại: Cao")

$p = $d.Paragraphs.Item(161)
$full = $p.Range
$full.Font.Italic = 0
$full.Font.ItalicBi = 0
$body = $p.Range
$body.MoveEnd(1, -1) | Out-Null
$body.Delete()
$body.InsertAfter("Giải pháp xử lý: Cải tiến chương trình hoặc do người dùng chưa cài đủ thư viện")

Write-Output "done"
